$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update recomputed cosinor statistics for rows 2-5 ---
# Row 2
$ws.Range("E2").Value = 25.50000000000055
$ws.Range("G2").Value = [double]"9.927504263096409e-09"
$ws.Range("H2").Value = [double]"3.172644230120319e-08"
$ws.Range("I2").Value = 0.6710926719905884
$ws.Range("K2").Value = 0.7954155341459099
$ws.Range("L2").Value = '[0.46030109520207674, 1.1305299730897431]'
$ws.Range("M2").Value = [double]"7.139653652066258e-06"
$ws.Range("N2").Value = [double]"7.139653652066258e-06"
$ws.Range("O2").Value = 2.446605690183042
$ws.Range("P2").Value = '[2.0440793041375027, 2.8491320762285817]'
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 36.60522751132503
$ws.Range("T2").Value = '[36.43306353592746, 36.77739148672261]'
$ws.Range("W2").Value = 15.5705705705709
$ws.Range("X2").Value = 13.93693693693723
$ws.Range("Y2").Value = 17.20420420420458
# Row 3
$ws.Range("C3").Value = '2_induction_dd'
$ws.Range("E3").Value = 23.86000000000029
$ws.Range("H3").Value = [double]"6.579099405186112e-16"
$ws.Range("K3").Value = 0.9995582072287164
$ws.Range("L3").Value = '[0.9285580815374104, 1.0705583329200223]'
$ws.Range("O3").Value = 2.19502669890458
$ws.Range("P3").Value = '[2.1195530015210418, 2.270500396288118]'
$ws.Range("S3").Value = 36.5335726112934
$ws.Range("T3").Value = '[36.49665871083259, 36.57048651175421]'
$ws.Range("W3").Value = 15.52452452452471
$ws.Range("X3").Value = 15.23791791791811
$ws.Range("Y3").Value = 15.81113113113132
# Row 4
$ws.Range("C4").Value = '3_hypo_dd'
$ws.Range("E4").Value = 23.82000000000028
$ws.Range("H4").Value = [double]"6.579099405186112e-16"
$ws.Range("K4").Value = 0.8244727870485775
$ws.Range("L4").Value = '[0.7355732654067211, 0.913372308690434]'
$ws.Range("O4").Value = 2.018921405009657
$ws.Range("P4").Value = '[1.9182898084982725, 2.119553001521041]'
$ws.Range("S4").Value = 36.39408304552244
$ws.Range("T4").Value = '[36.346749509543784, 36.44141658150109]'
$ws.Range("W4").Value = 16.16612612612632
$ws.Range("X4").Value = 15.78462462462482
$ws.Range("Y4").Value = 16.54762762762783
# Row 5
$ws.Range("C5").Value = '4_hypo_dl'
$ws.Range("E5").Value = 24.28000000000036
$ws.Range("G5").Value = [double]"7.294165271787278e-14"
$ws.Range("H5").Value = [double]"4.103569980914167e-13"
$ws.Range("K5").Value = 0.8174841104504722
$ws.Range("L5").Value = '[0.5558681172994326, 1.0791001036015118]'
$ws.Range("M5").Value = [double]"3.164085882190193e-09"
$ws.Range("N5").Value = [double]"4.21878117625359e-09"
$ws.Range("O5").Value = -0.779894872963232
$ws.Range("P5").Value = '[-1.0943686120613094, -0.4654211338651546]'
$ws.Range("Q5").Value = [double]"1.898453402704448e-06"
$ws.Range("R5").Value = [double]"1.898453402704448e-06"
$ws.Range("S5").Value = 36.26295110014394
$ws.Range("T5").Value = '[36.125057529025696, 36.40084467126218]'
$ws.Range("W5").Value = 3.013733733733776
$ws.Range("X5").Value = 1.798518518518543
$ws.Range("Y5").Value = 4.228948948949009

# --- Blank cells that must remain empty-string (Text) typed, not fully removed ---
$ws.Range("I3").Value = "'"
$ws.Range("I3").Style = "Normal"
$ws.Range("I5").Value = "'"
$ws.Range("I5").Style = "Normal"

# --- Remove the now-obsolete 6th row (5th animal/test entry dropped) ---
$ws.Rows.Item(6).Delete()
